$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new inventory row (row 77) following the existing pattern.
$ws.Range("A77").Value = "9MXYL2"
$ws.Range("B77").Value = "Luz led romantica"
$ws.Range("D77").Value = 30000
$ws.Range("E77").Value = 130000
$ws.Range("F77").Value = 4
$ws.Range("G77").Value = 4
$ws.Range("H77").Formula = "=(E77-D77)*G77"
$ws.Range("I77").Formula = "=D77*F77"
$ws.Range("J77").Value = 120000
